$d = $word.ActiveDocument

# Locate the two list-item paragraphs involved in the merge:
#   pRegulation : "Régulation des moteurs selon les données de l'IMU"
#   pCommande   : "Commande des moteurs selon la régulation et les capteurs de
#                  proximité" + " (« rebonds »)"
$pRegulation = $null
$pCommande = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "R*gulation des moteurs selon les donn*es de l*IMU*") {
        $pRegulation = $p
    } elseif ($t -like "Commande des moteurs selon la r*gulation*") {
        $pCommande = $p
    }
}

if ($pRegulation -eq $null -or $pCommande -eq $null) {
    throw "Could not locate the target paragraphs"
}

# Remove the whole "Commande ..." paragraph (text + its own paragraph mark) so
# the "Régulation ..." paragraph becomes directly followed by "Main :" again,
# while keeping the "Régulation ..." paragraph's own formatting (pStyle /
# numPr) untouched.
$pCommande.Range.Delete()

# Replace the "Régulation ..." paragraph's text (but not its trailing
# paragraph mark) with the new three-run sentence describing both manoeuvres:
#   "Commande" + " des moteurs selon les données de l'IMU" +
#   " & des capteurs de proximité"
$target = $d.Range($pRegulation.Range.Start, $pRegulation.Range.End - 1)

$eacute = [char]0x00E9
$rsquo = [char]0x2019

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>Commande</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> des moteurs selon les donn' + $eacute + 'es de l' + $rsquo + 'IMU</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> &amp; des capteurs de proximit' + $eacute + '</w:t></w:r>' + `
    '</w:p>'

$target.InsertXML($newXml)
